$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values look like plain numbers (single decimal point) and
# would otherwise be auto-converted from text to a numeric type when
# assigned. Force those specific cells to keep a Text number format so
# the stored value matches the original inline-string (text) cell type.
$numericLookingCells = @("D5","D10","D11","D19","D20","D23","D25","D27","D36","D39","D40","D42","D44","D45","D48","D49","D51")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "27.001.22"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.561.24"
$ws.Range("E3").Value = "  +0.57%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.09%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "207.30"
$ws.Range("E5").Value = "  +0.28%  "

# Row 6 (XRP)
$ws.Range("E6").Value = "  +1.23%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  -0.11%  "

# Row 8 (Solana)
$ws.Range("E8").Value = "  +2.10%  "

# Row 9 (Cardano)
$ws.Range("E9").Value = "  +0.06%  "

# Row 10 (Dogecoin)
$ws.Range("D10").Value = "0.0596"
$ws.Range("E10").Value = "  +1.84%  "

# Row 11 (TRON)
$ws.Range("D11").Value = "0.0859"
$ws.Range("E11").Value = "  +0.18%  "

# Row 12 (WrappedliquidstakedEther2.0)
$ws.Range("D12").Value = "1.784.09"
$ws.Range("E12").Value = "  +0.62%  "

# Row 13 (WrappedEther)
$ws.Range("D13").Value = "1.561.97"
$ws.Range("E13").Value = "  +0.61%  "

# Row 14 (Polkadot)
$ws.Range("E14").Value = "  +1.01%  "

# Row 15 (Polygon)
$ws.Range("E15").Value = "  +1.31%  "

# Row 16 (Litecoin)
$ws.Range("E16").Value = "  +0.75%  "

# Row 17 (WrappedBTC)
$ws.Range("D17").Value = "27.000.50"
$ws.Range("E17").Value = "  +0.46%  "

# Row 18 (ShibaInu)
$ws.Range("D18").Value = "0.0₃0705"
$ws.Range("E18").Value = "  +2.38%  "

# Row 19 (BitcoinCash)
$ws.Range("D19").Value = "217.08"
$ws.Range("E19").Value = "  +0.16%  "

# Row 20 (Chainlink)
$ws.Range("D20").Value = "7.36"
$ws.Range("E20").Value = "  +2.19%  "

# Row 21 (Dai)
$ws.Range("E21").Value = "  -0.09%  "

# Row 22 (Uniswap)
$ws.Range("E22").Value = "  +1.67%  "

# Row 23 (Avalanche)
$ws.Range("D23").Value = "9.26"
$ws.Range("E23").Value = "  +0.94%  "

# Row 24 (Toncoin)
$ws.Range("E24").Value = "  -2.34%  "

# Row 25 (Monero)
$ws.Range("D25").Value = "153.18"
$ws.Range("E25").Value = "  -0.17%  "

# Row 26 (Cosmos)
$ws.Range("E26").Value = "  +0.08%  "

# Row 27 (EthereumClassic)
$ws.Range("D27").Value = "15.06"
$ws.Range("E27").Value = "  +1.29%  "

# Row 28 (Stellar)
$ws.Range("E28").Value = "  +1.40%  "

# Row 29 (BinanceUSD)
$ws.Range("E29").Value = "  -0.10%  "

# Row 30 (Hedera)
$ws.Range("E30").Value = "  +1.13%  "

# Row 31 (PancakeSwap)
$ws.Range("E31").Value = "  +1.96%  "

# Row 32 (Filecoin)
$ws.Range("E32").Value = "  +0.84%  "

# Row 33 (Maker)
$ws.Range("D33").Value = "1.424.38"
$ws.Range("E33").Value = "  +0.53%  "

# Row 34 (InternetComputer(DFINITY))
$ws.Range("E34").Value = "  +3.71%  "

# Row 35 (LidoDAOToken)
$ws.Range("E35").Value = "  +3.50%  "

# Row 36 (TrustWalletToken)
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +9.20%  "

# Row 37 (HuobiToken)
$ws.Range("E37").Value = "  +1.34%  "

# Row 38 (VeChain)
$ws.Range("E38").Value = "  +0.83%  "

# Row 39 (ImmutableX)
$ws.Range("D39").Value = "0.533"
$ws.Range("E39").Value = "  +2.08%  "

# Row 40 (ARBITRUM)
$ws.Range("D40").Value = "0.809"
$ws.Range("E40").Value = "  +0.38%  "

# Row 41 (PaxDollar)
$ws.Range("E41").Value = "  -0.13%  "

# Row 42 (FraxShare)
$ws.Range("D42").Value = "5.69"
$ws.Range("E42").Value = "  +1.04%  "

# Row 43 (MXToken)
$ws.Range("E43").Value = "  +2.87%  "

# Row 44 (WEMIXToken)
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +1.15%  "

# Row 45 (Aave)
$ws.Range("D45").Value = "64.97"
$ws.Range("E45").Value = "  +2.12%  "

# Row 46 (RenderToken)
$ws.Range("E46").Value = "  +0.30%  "

# Row 47 (RocketPoolETH)
$ws.Range("D47").Value = "1.698.02"
$ws.Range("E47").Value = "  +0.63%  "

# Row 48 (Quant) - only price changes, Volume(1h) unchanged
$ws.Range("D48").Value = "87.53"

# Row 49 (Cronos)
$ws.Range("D49").Value = "0.0521"
$ws.Range("E49").Value = "  +0.66%  "

# Row 50 (BabyDogeCoin)
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  +0.24%  "

# Row 51 (Algorand)
$ws.Range("D51").Value = "0.0955"
$ws.Range("E51").Value = "  -0.45%  "
